$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 697.8
$ws.Range("I12").Value = 697.8
$ws.Range("K12").Value = 697.8
$ws.Range("M12").Value = -527.8

$ws.Range("H98").Value = 1076.72
$ws.Range("I98").Value = 823.05554
$ws.Range("J98").Value = 1729
$ws.Range("K98").Value = 823.05554
$ws.Range("L98").Value = 1729
$ws.Range("M98").Value = 674.94446
$ws.Range("N98").Value = -4725

$ws.Range("H101").Value = 1146.8334
$ws.Range("J101").Value = 1226.8572
$ws.Range("L101").Value = 3680.5716
$ws.Range("N101").Value = -6924.571599999999

$ws.Range("H122").Value = 1076.72
$ws.Range("I122").Value = 823.05554
$ws.Range("J122").Value = 1729
$ws.Range("K122").Value = 2469.16662
$ws.Range("L122").Value = 5187
$ws.Range("M122").Value = -19.16661999999997
$ws.Range("N122").Value = -10087

$ws.Range("H138").Value = 3364.6
$ws.Range("I138").Value = 3344.5
$ws.Range("J138").Value = 3378
$ws.Range("K138").Value = 10033.5
$ws.Range("L138").Value = 10134
$ws.Range("M138").Value = -4893.5
$ws.Range("N138").Value = -20414

$ws.Range("H141").Value = 6971.1333
$ws.Range("I141").Value = 7333.075
$ws.Range("J141").Value = 4075.6
$ws.Range("K141").Value = 21999.225
$ws.Range("L141").Value = 12226.8
$ws.Range("M141").Value = -16819.225
$ws.Range("N141").Value = -22586.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4275.84
$ws.Range("I32").Value = 3522.2842
$ws.Range("J32").Value = 9801.916999999999
$ws.Range("K32").Value = 3522.2842
$ws.Range("L32").Value = 9801.916999999999
$ws.Range("M32").Value = -3235.2842
$ws.Range("N32").Value = -10375.917

$ws.Range("H74").Value = 1505.7931
$ws.Range("I74").Value = 1505.7931
$ws.Range("K74").Value = 1505.7931
$ws.Range("M74").Value = -631.7931000000001

$ws.Range("H77").Value = 1505.7931
$ws.Range("I77").Value = 1505.7931
$ws.Range("K77").Value = 7528.9655
$ws.Range("M77").Value = -3160.9655

$ws.Range("H132").Value = 2040.08
$ws.Range("I132").Value = 2231.2715
$ws.Range("J132").Value = 1593.9667
$ws.Range("K132").Value = 6693.814499999999
$ws.Range("L132").Value = 4781.9001
$ws.Range("M132").Value = -4163.814499999999
$ws.Range("N132").Value = -9841.900099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 434.44446
$ws.Range("I22").Value = 236
$ws.Range("J22").Value = 593.2
$ws.Range("K22").Value = 236
$ws.Range("L22").Value = 593.2
$ws.Range("M22").Value = -63
$ws.Range("N22").Value = -939.2

$ws.Range("H134").Value = 2845.6956
$ws.Range("I134").Value = 2318.0513
$ws.Range("J134").Value = 5785.4287
$ws.Range("K134").Value = 6954.1539
$ws.Range("L134").Value = 17356.2861
$ws.Range("M134").Value = -4419.1539
$ws.Range("N134").Value = -22426.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7824.5625
$ws.Range("I22").Value = 23849.2
$ws.Range("J22").Value = 540.63635
$ws.Range("K22").Value = 23849.2
$ws.Range("L22").Value = 540.63635
$ws.Range("M22").Value = -23499.2
$ws.Range("N22").Value = -1240.63635

$ws.Range("H31").Value = 3812.9722
$ws.Range("I31").Value = 2530.2083
$ws.Range("J31").Value = 6378.5
$ws.Range("K31").Value = 2530.2083
$ws.Range("L31").Value = 6378.5
$ws.Range("M31").Value = -2235.2083
$ws.Range("N31").Value = -6968.5

$ws.Range("H34").Value = 3812.9722
$ws.Range("I34").Value = 2530.2083
$ws.Range("J34").Value = 6378.5
$ws.Range("K34").Value = 2530.2083
$ws.Range("L34").Value = 6378.5
$ws.Range("M34").Value = -2328.2083
$ws.Range("N34").Value = -6782.5

$ws.Range("H58").Value = 2200.5625
$ws.Range("I58").Value = 1781.3636
$ws.Range("J58").Value = 3122.8
$ws.Range("K58").Value = 1781.3636
$ws.Range("L58").Value = 3122.8
$ws.Range("M58").Value = -1578.3636
$ws.Range("N58").Value = -3528.8

$ws.Range("H132").Value = 5391.712
$ws.Range("I132").Value = 1387.909
$ws.Range("J132").Value = 17136.2
$ws.Range("K132").Value = 4163.727000000001
$ws.Range("L132").Value = 51408.60000000001
$ws.Range("M132").Value = -1633.727000000001
$ws.Range("N132").Value = -56468.60000000001

$ws.Range("H134").Value = 2981.592
$ws.Range("I134").Value = 2777.8057
$ws.Range("J134").Value = 6649.75
$ws.Range("K134").Value = 8333.417099999999
$ws.Range("L134").Value = 19949.25
$ws.Range("M134").Value = -5798.417099999999
$ws.Range("N134").Value = -25019.25

$ws.Range("H136").Value = 2200.5625
$ws.Range("I136").Value = 1781.3636
$ws.Range("J136").Value = 3122.8
$ws.Range("K136").Value = 5344.0908
$ws.Range("L136").Value = 9368.400000000001
$ws.Range("M136").Value = -2794.0908
$ws.Range("N136").Value = -14468.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3987.6667
$ws.Range("I87").Value = 3861.125
$ws.Range("K87").Value = 11583.375
$ws.Range("M87").Value = -10335.375

$ws.Range("H90").Value = 3987.6667
$ws.Range("I90").Value = 3861.125
$ws.Range("K90").Value = 34750.125
$ws.Range("M90").Value = -28510.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4389.1
$ws.Range("I31").Value = 5518
$ws.Range("J31").Value = 3905.2856
$ws.Range("K31").Value = 5518
$ws.Range("L31").Value = 3905.2856
$ws.Range("M31").Value = -5270
$ws.Range("N31").Value = -4401.2856

$ws.Range("H136").Value = 6631.9546
$ws.Range("I136").Value = 6919.0586
$ws.Range("J136").Value = 5655.8
$ws.Range("K136").Value = 20757.1758
$ws.Range("L136").Value = 16967.4
$ws.Range("M136").Value = -18207.1758
$ws.Range("N136").Value = -22067.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 671.41174
$ws.Range("I100").Value = 654.3333
$ws.Range("K100").Value = 1308.6666
$ws.Range("M100").Value = -767.6666

$ws.Range("H132").Value = 3304.848
$ws.Range("I132").Value = 2861.186
$ws.Range("J132").Value = 9664
$ws.Range("K132").Value = 8583.558000000001
$ws.Range("L132").Value = 28992
$ws.Range("M132").Value = -6053.558000000001
$ws.Range("N132").Value = -34052

$ws.Range("H136").Value = 3971.1428
$ws.Range("I136").Value = 3077.7222
$ws.Range("J136").Value = 9331.666999999999
$ws.Range("K136").Value = 9233.1666
$ws.Range("L136").Value = 27995.001
$ws.Range("M136").Value = -6683.1666
$ws.Range("N136").Value = -33095.001

